$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextCell 'D2' '26.370.20'
Set-TextCell 'E2' '  +3.54%  '
Set-TextCell 'D3' '1.724.37'
Set-TextCell 'E3' '  +3.28%  '
Set-TextCell 'E4' '  -0.04%  '
Set-TextCell 'D5' '242.66'
Set-TextCell 'E5' '  +2.12%  '
Set-TextCell 'D6' '0.9994'
Set-TextCell 'E6' '  -0.09%  '
Set-TextCell 'D7' '0.4742'
Set-TextCell 'E7' '  -1.12%  '
Set-TextCell 'D8' '0.2642'
Set-TextCell 'E8' '  +0.42%  '
Set-TextCell 'D9' '0.06206'
Set-TextCell 'E9' '  +0.45%  '
Set-TextCell 'D10' '1.719.84'
Set-TextCell 'E10' '  +3.04%  '
Set-TextCell 'D11' '0.07074'
Set-TextCell 'E11' '  +1.00%  '
Set-TextCell 'D12' '15.52'
Set-TextCell 'E12' '  +4.54%  '
Set-TextCell 'D13' '0.5982'
Set-TextCell 'E13' '  +1.60%  '
Set-TextCell 'D14' '4.432'
Set-TextCell 'E14' '  +1.41%  '
Set-TextCell 'D15' '76.45'
Set-TextCell 'E15' '  +1.98%  '
Set-TextCell 'D16' '0.9994'
Set-TextCell 'E16' '  -0.09%  '
Set-TextCell 'D17' '26.375.85'
Set-TextCell 'E17' '  +3.57%  '
Set-TextCell 'E18' '  -0.07%  '
Set-TextCell 'D19' '0.000006838'
Set-TextCell 'E19' '  +1.16%  '
Set-TextCell 'D20' '11.58'
Set-TextCell 'E20' '  +1.08%  '
Set-TextCell 'D21' '1.939.81'
Set-TextCell 'E21' '  +3.15%  '
Set-TextCell 'D22' '4.524'
Set-TextCell 'E22' '  +1.49%  '
Set-TextCell 'D23' '8.762'
Set-TextCell 'E23' '  +0.28%  '
Set-TextCell 'D24' '5.266'
Set-TextCell 'E24' '  -0.40%  '
Set-TextCell 'D25' '135.32'
Set-TextCell 'E25' '  -1.31%  '
Set-TextCell 'D26' '15.23'
Set-TextCell 'E26' '  +1.39%  '
Set-TextCell 'D27' '1.773'
Set-TextCell 'E27' '  +2.94%  '
Set-TextCell 'D28' '1.401'
Set-TextCell 'E28' '  +0.81%  '
Set-TextCell 'D29' '106.93'
Set-TextCell 'E29' '  +1.99%  '
Set-TextCell 'D30' '3.958'
Set-TextCell 'E30' '  +0.32%  '
Set-TextCell 'D31' '3.694'
Set-TextCell 'E31' '  +1.35%  '
Set-TextCell 'D32' '0.07809'
Set-TextCell 'E32' '  +0.02%  '
Set-TextCell 'D33' '0.04508'
Set-TextCell 'E33' '  +6.71%  '
Set-TextCell 'B34' 'HuobiToken'
Set-TextCell 'C34' 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
Set-TextCell 'D34' '2.616'
Set-TextCell 'E34' '  +0.44%  '
Set-TextCell 'B35' 'ARBITRUM'
Set-TextCell 'C35' 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextCell 'D35' '0.9829'
Set-TextCell 'E35' '  +3.62%  '
Set-TextCell 'B36' 'ImmutableX'
Set-TextCell 'C36' 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextCell 'D36' '0.6234'
Set-TextCell 'E36' '  +2.33%  '
Set-TextCell 'B37' 'TrustWalletToken'
Set-TextCell 'C37' 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextCell 'D37' '0.9457'
Set-TextCell 'E37' '  +10.46%  '
Set-TextCell 'B38' 'Quant'
Set-TextCell 'C38' 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
Set-TextCell 'D38' '114.42'
Set-TextCell 'E38' '  +18.97%  '
Set-TextCell 'B39' 'MXToken'
Set-TextCell 'C39' 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextCell 'D39' '2.456'
Set-TextCell 'E39' '  -5.52%  '
Set-TextCell 'B40' 'RenderToken'
Set-TextCell 'C40' 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextCell 'D40' '1.936'
Set-TextCell 'E40' '  +4.61%  '
Set-TextCell 'B41' 'PaxDollar'
Set-TextCell 'C41' 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
Set-TextCell 'D41' '0.9999'
Set-TextCell 'E41' '  +0.03%  '
Set-TextCell 'B42' 'FraxShare'
Set-TextCell 'C42' 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextCell 'D42' '5.696'
Set-TextCell 'E42' '  +17.65%  '
Set-TextCell 'B43' 'VeChain'
Set-TextCell 'C43' 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextCell 'D43' '0.01487'
Set-TextCell 'E43' '  +0.67%  '
Set-TextCell 'B44' 'TheSandbox'
Set-TextCell 'C44' 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
Set-TextCell 'D44' '0.3837'
Set-TextCell 'E44' '  +1.79%  '
Set-TextCell 'B45' 'Algorand'
Set-TextCell 'C45' 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-TextCell 'D45' '0.1188'
Set-TextCell 'E45' '  +6.15%  '
Set-TextCell 'B46' 'Aptos'
Set-TextCell 'C46' 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextCell 'D46' '6.383'
Set-TextCell 'E46' '  +2.80%  '
Set-TextCell 'B47' 'Cronos'
Set-TextCell 'C47' 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextCell 'D47' '0.05274'
Set-TextCell 'E47' '  +0.45%  '
Set-TextCell 'B48' 'EnergySwap'
Set-TextCell 'C48' 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextCell 'D48' '7.868'
Set-TextCell 'E48' '  +6.63%  '
Set-TextCell 'B49' 'Elrond'
Set-TextCell 'C49' 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
Set-TextCell 'D49' '30.43'
Set-TextCell 'E49' '  +2.09%  '
Set-TextCell 'B50' 'Decentraland'
Set-TextCell 'C50' 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
Set-TextCell 'D50' '0.3391'
Set-TextCell 'E50' '  +1.79%  '
Set-TextCell 'B51' 'NEARProtocol'
Set-TextCell 'C51' 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextCell 'D51' '1.219'
Set-TextCell 'E51' '  +2.04%  '
